# COP21_man_bosnia_herz.xlsx update ("Updated documents, Local folder file switch")
#
# Content edits:
#  - Row 4  (H4):  revise the 30-word explanation text
#  - Row 6  (D6:I6): Unit/Scale/Time/Principle/30-word-explanation/Notes refresh
#  - Row 11 (G11/H11): Principle/30-word-explanation refresh
#
# View state: the saved file now has the viewport/selection resting on D6
# instead of H11, so move the active selection there to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 4 - 30 word explanation ---
$ws.Range("H4").Value = "Prescribing the need to take on action for the benefit of all. "

# --- Row 6 - Unit, Scale, Time, Principle, 30 word explanation, Notes ---
$ws.Range("D6").Value = "n.a."
$ws.Range("E6").Value = "regional"
$ws.Range("F6").Value = "n.a."
$ws.Range("G6").Value = "prioritarian"
$ws.Range("H6").Value = "Indicating the need for support presented by a developed country to help worse off. "
$ws.Range("I6").Value = "allign with EU"

# --- Row 11 - Principle, 30 word explanation ---
$ws.Range("G11").Value = "egalitarian"
$ws.Range("H11").Value = "Prescribing the need for fair policy with clear distrinction and defined relations, mentioning shared responsilibity. "

# --- Row heights settle ~4% shorter after the resave (wrapped-text rows
#     recompute against the new default font metrics) ---
$ws.Rows.Item(2).RowHeight = 86.4
$ws.Rows.Item(3).RowHeight = 115.2
$ws.Rows.Item(4).RowHeight = 72
$ws.Rows.Item(5).RowHeight = 158.4
$ws.Rows.Item(6).RowHeight = 86.4
$ws.Rows.Item(7).RowHeight = 172.8
$ws.Rows.Item(8).RowHeight = 115.2
$ws.Rows.Item(9).RowHeight = 129.6
$ws.Rows.Item(10).RowHeight = 144
$ws.Rows.Item(11).RowHeight = 201.6
$ws.Rows.Item(12).RowHeight = 172.8
$ws.Rows.Item(13).RowHeight = 57.6

# --- Move selection/active cell to D6 (matches the saved view state) ---
$ws.Range("D6").Select()
